# Update the "想去人数" (want-to-go count) values in column F
# on the "展览" and "全部类型" sheets, reflecting newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (rows 1-38) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 446
$ws1.Range("F5").Value = 1753
$ws1.Range("F7").Value = 2218
$ws1.Range("F14").Value = 310
$ws1.Range("F17").Value = 191
$ws1.Range("F21").Value = 3980
$ws1.Range("F24").Value = 29
$ws1.Range("F26").Value = 113
$ws1.Range("F28").Value = 27
$ws1.Range("F30").Value = 94
$ws1.Range("F34").Value = 1017
$ws1.Range("F36").Value = 2582
$ws1.Range("F37").Value = 431
$ws1.Range("F38").Value = 26

# --- Sheet "全部类型" (rows 1-39) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 446
$ws4.Range("F5").Value = 1753
$ws4.Range("F7").Value = 2218
$ws4.Range("F14").Value = 310
$ws4.Range("F17").Value = 191
$ws4.Range("F21").Value = 3980
$ws4.Range("F24").Value = 29
$ws4.Range("F26").Value = 113
$ws4.Range("F28").Value = 27
$ws4.Range("F30").Value = 95
$ws4.Range("F35").Value = 1017
$ws4.Range("F37").Value = 2582
$ws4.Range("F38").Value = 431
$ws4.Range("F39").Value = 26
